$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("/NewDataSet/Table[2]/Location", "[A-Z a-z].*"),
    @("/NewDataSet/Table[2]/PostCode", "[A-Z a-z 0-9].*"),
    @("/NewDataSet/Table[3]/Location", "[A-Z a-z].*"),
    @("/NewDataSet/Table[3]/PostCode", "[A-Z a-z 0-9].*"),
    @("/NewDataSet/Table[4]/Location", "[A-Z a-z].*"),
    @("/NewDataSet/Table[4]/PostCode", "[A-Z a-z 0-9].*")
)

$row = 4
foreach ($pair in $data) {
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
    $row++
}

$ws.Range("B4:C9").Select()
